$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the "Location" values in column D (shared-string driven) ---
# ERPLZI -> 16_ERPLZI_2km, ERPLS2 -> 10_ERPLS2_3km, DAPLD -> 10_DAPLD_6km
$col = $ws.Range("D1:D66")
$col.Replace("ERPLZI", "16_ERPLZI_2km")
$col.Replace("ERPLS2", "10_ERPLS2_3km")
$col.Replace("DAPLD", "10_DAPLD_6km")

# --- Widen column D and drop the old auto "best fit" sizing ---
$ws.Columns("D").ColumnWidth = 24.857142857142858

# --- Turn on AutoFilter for the table range ---
$ws.Range("A1:F66").AutoFilter()

# --- Register the hidden _FilterDatabase defined name AutoFilter creates ---
$fdb = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$F`$66")
$fdb.Visible = $false

# --- Move the active selection ---
$ws.Range("D14").Select()
